$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 86, pushing the existing rows
# (old rows 86-181) down to rows 87-182.
$ws.Rows.Item(86).Insert()

# Populate the new row 86 with the new weekly price entry.
$ws.Range("A86").Value = 4
$ws.Range("B86").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C86").Value = "Los Lagos"
$ws.Range("D86").Value = 44483
$ws.Range("E86").Value = 10
$ws.Range("F86").Value = 100112045
$ws.Range("G86").Value = "Zapallo"
$ws.Range("H86").Value = "Paine"
$ws.Range("I86").Value = "1a (guarda)"
$ws.Range("J86").Value = 500
$ws.Range("K86").Value = 480
$ws.Range("L86").Value = 480
$ws.Range("M86").Value = 480
$ws.Range("N86").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O86").Value = "Región Metropolitana"
$ws.Range("P86").Value = 480
$ws.Range("Q86").Value = 1
$ws.Range("R86").Value = "Hortaliza"
